$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10, shifting existing row 10 (and below) down by one.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 44505
$ws.Cells.Item(10, 4).Style = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100112022
$ws.Cells.Item(10, 7).Value = "Arveja Verde"
$ws.Cells.Item(10, 8).Value = "Perfection"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 210
$ws.Cells.Item(10, 11).Value = 6500
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 6714
$ws.Cells.Item(10, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 269
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
